$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L header + values ("not 1988 + entire_scale")
$ws.Range("L5").Value = "not 1988 + entire_scale"

# Updated R^2_Train row (row 6)
$ws.Range("F6").Value = 0.95713199999999998
$ws.Range("G6").Value = 0.95050230000000002
$ws.Range("H6").Value = 0.95459760000000005
$ws.Range("I6").Value = 0.95386249999999995
$ws.Range("J6").Value = 0.95804310000000004
$ws.Range("L6").Value = 0.90412029999999999

# Updated R^2_Test row (row 7)
$ws.Range("F7").Value = 0.66669560000000005
$ws.Range("G7").Value = 0.67541430000000002
$ws.Range("H7").Value = 0.69006259999999997
$ws.Range("I7").Value = 0.68827439999999995
$ws.Range("J7").Value = 0.71332530000000005
$ws.Range("L7").Value = 0.58579820000000005

# Updated RSME_Train row (row 8)
$ws.Range("F8").Value = 0.72757309999999997
$ws.Range("G8").Value = 0.74498299999999995
$ws.Range("H8").Value = 0.71458350000000004
$ws.Range("I8").Value = 0.7199875
$ws.Range("J8").Value = 0.69620930000000003
$ws.Range("L8").Value = 1.074109

# Updated RSME_Test row (row 9)
$ws.Range("F9").Value = 1.9058470000000001
$ws.Range("G9").Value = 1.8668720000000001
$ws.Range("H9").Value = 1.8236129999999999
$ws.Range("I9").Value = 1.8267960000000001
$ws.Range("J9").Value = 1.757136
$ws.Range("L9").Value = 2.1296249999999999

# Update selection / view state to match the saved workbook
$ws.Range("J9").Select()
